$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.867.37"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.18%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.643.02"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.78%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'215.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.18%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.5062"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.03%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.20%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2591"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.90%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06439"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.65%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'20.47"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +5.13%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07815"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.88%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.35%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'1.640.20"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.60%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'1.868.73"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.68%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.5639"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.89%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.0₅7700"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.65%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'63.46"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.29%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'25.877.48"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.14%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'1.003"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.11%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'194.36"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.04%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'4.377"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.94%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'9.966"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.94%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'6.127"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.33%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'1.003"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.13%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'1.799"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -6.08%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'140.48"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.01%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.1247"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.90%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'6.846"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.89%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'15.56"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.05%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.246"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.59%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'0.04953"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.57%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'3.322"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.57%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'3.257"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +2.44%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'1.588"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +3.07%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.382"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.35%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.9098"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.62%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.579"
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Value = "'0.5577"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.06%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'1.132.06"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.33%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.01567"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.94%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.14%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'5.547"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.36%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.8048"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.83%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'98.57"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.54%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'1.777.34"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.50%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.0₈113"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -4.88%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'55.88"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +2.30%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.86%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'7.803"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +3.97%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.06%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'1.0000"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.29%  "
$ws.Range('E51').Style = 'Normal'
